# Mass Interview uncheck for Only Available Interviewers
#
# Appends the latest sprint-run rows to the AMSIN, BETA and AMS history
# sheets (INDIA is untouched this round) and normalises the formatting of
# the previously "newest" row on AMSIN (row 23) now that it is no longer
# the last one in the file.
#
# Notes on technique (all driven purely through the Excel COM object
# model):
#  - Date-looking strings ("2021-10-06", ...) must stay literal text
#    (matching the existing "Run Date" column), so instead of assigning
#    them straight to .Value (which Excel would parse as a date serial)
#    we build the text with a formula ( ="2021-10-06" ), then convert
#    that formula to a literal value in place via Copy/PasteSpecial
#    (values only). This keeps the natural column formatting.
#  - Brand new cells pick up their column's default style automatically
#    the first time they're written, and ClearContents() on an existing
#    cell (keeping the cell but dropping its content) makes the next
#    write behave like a "new" cell too - this is how row 23 on AMSIN
#    picks up the same style as the rest of its column.
#  - Column B (Run Time) needs the workbook's custom date/time style
#    already used elsewhere in the same column, so after writing the
#    numeric value we copy that formatting over with PasteSpecial
#    (formats only) from a neighbouring cell that already has it.
#  - AMS row 18 is the most recently appended row overall, so (matching
#    the pattern already present in the workbook) it is intentionally
#    left without the normalised column style on columns A and C..G.
#    Cells that are new (beyond the sheet's previous used range) pick up
#    the column's style as soon as they're written, so those values are
#    staged in a scratch row just below, then moved (Cut) into place -
#    a cut/move keeps the (unstyled) source formatting instead of
#    inheriting the destination column's default - and the scratch row
#    is deleted afterwards.

$wb = $excel.ActiveWorkbook

# ============================================================
# AMSIN: normalise row 23's styling, then append rows 24 and 25
# ============================================================
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# -- row 23 already holds the right data; re-enter it so it picks up
#    the same column-default style ("s=5") that every other row uses --
$wsAmsin.Range("A23:G23").ClearContents()

$wsAmsin.Cells.Item(23, 1).Formula = '="2021-10-06"'
$wsAmsin.Range("A23").Copy()
$wsAmsin.Range("A23").PasteSpecial(-4163)

$wsAmsin.Cells.Item(23, 2).Value = 44475.70605298611
$wsAmsin.Cells.Item(23, 3).Value = "151_regression"
$wsAmsin.Cells.Item(23, 4).Value = 75
$wsAmsin.Cells.Item(23, 5).Value = 73
$wsAmsin.Cells.Item(23, 6).Value = 2
$wsAmsin.Cells.Item(23, 7).Value = 3.15

# -- append row 24 (151 -> 152 first-cycle run) --
$wsAmsin.Cells.Item(24, 1).Formula = '="2021-10-26"'
$wsAmsin.Range("A24").Copy()
$wsAmsin.Range("A24").PasteSpecial(-4163)

$wsAmsin.Cells.Item(24, 2).Value = 44495.65657489583
$wsAmsin.Range("B23").Copy()
$wsAmsin.Range("B24").PasteSpecial(-4122)

$wsAmsin.Cells.Item(24, 3).Value = "152_fstcycle"
$wsAmsin.Cells.Item(24, 4).Value = 75
$wsAmsin.Cells.Item(24, 5).Value = 73
$wsAmsin.Cells.Item(24, 6).Value = 2
$wsAmsin.Cells.Item(24, 7).Value = 2.51

# -- append row 25 (152 final regression run) --
$wsAmsin.Cells.Item(25, 1).Formula = '="2021-10-28"'
$wsAmsin.Range("A25").Copy()
$wsAmsin.Range("A25").PasteSpecial(-4163)

$wsAmsin.Cells.Item(25, 2).Value = 44497.38869483797
$wsAmsin.Range("B23").Copy()
$wsAmsin.Range("B25").PasteSpecial(-4122)

$wsAmsin.Cells.Item(25, 3).Value = "152_fnlrgrsn"
$wsAmsin.Cells.Item(25, 4).Value = 75
$wsAmsin.Cells.Item(25, 5).Value = 73
$wsAmsin.Cells.Item(25, 6).Value = 2
$wsAmsin.Cells.Item(25, 7).Value = 2.32

# ============================================================
# BETA: append row 15 (152 beta run)
# ============================================================
$wsBeta = $wb.Worksheets.Item("BETA")

$wsBeta.Cells.Item(15, 1).Formula = '="2021-10-28"'
$wsBeta.Range("A15").Copy()
$wsBeta.Range("A15").PasteSpecial(-4163)

$wsBeta.Cells.Item(15, 2).Value = 44497.62782178241
$wsBeta.Range("B14").Copy()
$wsBeta.Range("B15").PasteSpecial(-4122)

$wsBeta.Cells.Item(15, 3).Value = "152_beta"
$wsBeta.Cells.Item(15, 4).Value = 75
$wsBeta.Cells.Item(15, 5).Value = 73
$wsBeta.Cells.Item(15, 6).Value = 2
$wsBeta.Cells.Item(15, 7).Value = 2.76

# ============================================================
# AMS: append row 18 (152 live test run) - the most recent run in
# this whole update, so it is left without the normalised style,
# same as AMSIN row 23 was before this update.
# ============================================================
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(19, 1).Formula = '="2021-10-28"'
$wsAms.Range("A19").Copy()
$wsAms.Range("A19").PasteSpecial(-4163)
$wsAms.Cells.Item(19, 1).Cut($wsAms.Cells.Item(18, 1))

$wsAms.Cells.Item(18, 2).Value = 44497.86153741145
$wsAms.Range("B17").Copy()
$wsAms.Range("B18").PasteSpecial(-4122)

$wsAms.Cells.Item(19, 3).Value = "152_livetest"
$wsAms.Cells.Item(19, 3).Cut($wsAms.Cells.Item(18, 3))

$wsAms.Cells.Item(19, 4).Value = 75
$wsAms.Cells.Item(19, 4).Cut($wsAms.Cells.Item(18, 4))

$wsAms.Cells.Item(19, 5).Value = 75
$wsAms.Cells.Item(19, 5).Cut($wsAms.Cells.Item(18, 5))

$wsAms.Cells.Item(19, 6).Value = 0
$wsAms.Cells.Item(19, 6).Cut($wsAms.Cells.Item(18, 6))

$wsAms.Cells.Item(19, 7).Value = 2.87
$wsAms.Cells.Item(19, 7).Cut($wsAms.Cells.Item(18, 7))

$wsAms.Rows.Item(19).Delete()
